$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Fill in the daily activity row (row 2)
$ws.Range("B2").Value = "Finish auto register script in Python"
$ws.Range("C2").Value = "Used Selenium, chrome driver, excel read and write, multiprocessing-Pool"

# Fill in the reflection / reminder row (row 6) and grow its height to fit the wrapped text
$ws.Range("B6").Value = "Reminder to modify"
$ws.Range("C6").Value = "1. make sure the input excel sheet with the proper name same as the village name;`n2. ban the register function before formal running, test the read and write Excel process;`n3. remove the user data after registration and avoid upload it github"
$ws.Rows.Item(6).RowHeight = 40.5

# Update the active selection to span the populated report area
$ws.Range("A1:C6").Select()
